$d = $word.ActiveDocument

$replacements = @(
    @{old = "698÷4="; new = "714÷4="},
    @{old = "819÷7="; new = "103÷8="},
    @{old = "802÷6="; new = "570÷5="},
    @{old = "649÷8="; new = "877÷9="},
    @{old = "284÷6="; new = "369÷2="},
    @{old = "623÷8="; new = "111÷2="},
    @{old = "609÷9="; new = "602÷9="},
    @{old = "848÷5="; new = "665÷5="},
    @{old = "872÷7="; new = "944÷7="},
    @{old = "490÷9="; new = "673÷6="},
    @{old = "116÷5="; new = "190÷9="},
    @{old = "191÷8="; new = "275÷9="},
    @{old = "860÷3="; new = "613÷7="},
    @{old = "840÷6="; new = "790÷7="},
    @{old = "202÷4="; new = "491÷6="},
    @{old = "703÷2="; new = "943÷7="},
    @{old = "648÷3="; new = "754÷7="},
    @{old = "384÷3="; new = "450÷9="},
    @{old = "751÷2="; new = "510÷3="},
    @{old = "812÷4="; new = "308÷9="},
    @{old = "298÷3="; new = "203÷7="},
    @{old = "310÷3="; new = "715÷4="},
    @{old = "343÷7="; new = "898÷3="},
    @{old = "933÷6="; new = "288÷7="},
    @{old = "243÷9="; new = "435÷5="}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
